$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend existing CHEBI-DRUG-ROLE entry (row 6) with additional drug role IDs
$ws.Range("D6").Value = '﻿Antidepressant drug [CHEBI:35469];Antimanic drug [CHEBI:35477];Antipsychotic drug [CHEBI:35476];Anxiolytic drug [CHEBI:35474];Central nervous system depressent [CHEBI:35488];Central nervous system drug [CHEBI:35470];Pharmaceutical [CHEBI:52217];Tranquilizing drug  [CHEBI:35473];Drug [CHEBI:23888]; psychotropic drug [CHEBI:35471]; opioid analgesic [CHEBI:35482]'

# Extend existing CHEBI-BIO-ROLE entry (row 5) with additional role IDs
$ws.Range("D5").Value = 'Carcinogen [CHEBI:50903];Alcohol dehydrogenase 1B [CHEBI:50269]; Flavouring agent [CHEBI:35617]; nicotinic acetylcholine receptor agonist [CHEBI:47958]; nicotinic antagonist [CHEBI:48878]; carcinogenic agent [CHEBI:50903]; cannabinoid receptor agonist [CHEBI:67072]; cb1 receptor antagonist [CHEBI:73416]; cb2 receptor agonist [CHEBI:73417]'

# Add new ERO ontology import row
$ws.Range("A26").Value = "ERO"
$ws.Range("B26").Value = "http://purl.obolibrary.org/obo/ero.owl"
$ws.Range("C26").Value = "process [BFO:0000015]"
$ws.Range("D26").Value = "intubation [ERO:0001108]"
$ws.Range("E26").Value = "minimal"

# Extend existing CHEBI-CHEM entry (row 4) with additional imported chemical IDs
$ws.Range("D4").Value = 'Amphetamine [CHEBI:2679];Baclofen [CHEBI:2972];Benzodiazepine [CHEBI:22720];Buprenorphine [CHEBI:3216];Bupropion [CHEBI:3219];Caffeine [CHEBI:27732];Cannabinoid [CHEBI:67194];Chemical substance [CHEBI:59999];Cocaine [CHEBI:27958];Codeine [CHEBI:38164];Ethanol [CHEBI:16236];Fentanyl [CHEBI:119915];Formaldehyde [CHEBI:16842];Gabapentin [CHEBI:42797];Heroin [CHEBI:27808];Methadone [CHEBI:6807];Methadrone  [CHEBI:59331];Methamphetamine [CHEBI:6809];Naloxone [CHEBI:7459];Nicotine [CHEBI:18723];4-(N-nitrosomethylamino)-1-(3-pyridyl)butan-1-one [CHEBI:32692];Synthetic cannabinoid [CHEBI:67201];Tetrahydrocannabinol [CHEBI:66964];Varenicline [CHEBI:84500];nitrosamine [CHEBI:35803]; Atomoxetine [CHEBI:127342];naphthalene [CHEBI:16482];propanal [CHEBI:17153];carbon monoxide [CHEBI:17245]; nicotine [CHEBI:18723]; carbonyl compound [CHEBI:36586]; elemental cadmium [CHEBI:37249]; pyrene [CHEBI:39106]; cytisine [CHEBI:4055]; crotonaldehyde [CHEBI:41607]; diazepam [CHEBI:49575]; hydrocodone [CHEBI:5779]; cannabidiol [CHEBI:69478]; oxycodone [CHEBI:7852]; rimonabant [CHEBI:34967]'

$ws.Range("D5").Select()
